$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "66.956.01"
$r.Style = "Normal"
$ws.Range("E2").Value = "  -2.04%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "2.632.58"
$r.Style = "Normal"
$ws.Range("E3").Value = "  -3.11%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "588.62"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -3.25%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "165.60"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  +0.01%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "0.541"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -2.22%  "
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = "2.631.77"
$r.Style = "Normal"
$ws.Range("E9").Value = "  -3.08%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.143"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -1.21%  "
$ws.Range("E11").Value = "  +1.36%  "
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "0.362"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -0.82%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "5.25"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -1.15%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "27.60"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -3.31%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "3.117.36"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -3.66%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "67.225.87"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -1.58%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "2.620.86"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -0.90%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "11.98"
$r.Style = "Normal"
$ws.Range("E19").Value = "  +0.83%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "8.08"
$r.Style = "Normal"
$ws.Range("E20").Value = "  +5.72%  "
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = "360.16"
$r.Style = "Normal"
$ws.Range("E21").Value = "  -3.10%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "4.34"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -3.56%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "4.74"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -4.62%  "
$ws.Range("E24").Value = "  +8.98%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "1.97"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -5.52%  "
$ws.Range("E26").Value = "  -0.04%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "70.66"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -3.27%  "
$ws.Range("E29").Value = "  +0.10%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "0.0000101"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -2.97%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "550.11"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -6.06%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "7.94"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -3.16%  "
$ws.Range("E33").Value = "  -4.13%  "
$ws.Range("E34").Value = "  -4.57%  "
$ws.Range("E35").Value = "  +2.40%  "
$ws.Range("E36").Value = "  +0.04%  "
$r = $ws.Range("D37")
$r.NumberFormat = "@"
$r.Value = "1.52"
$r.Style = "Normal"
$ws.Range("E37").Value = "  -4.80%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "157.47"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -3.32%  "
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = "19.18"
$r.Style = "Normal"
$ws.Range("E39").Value = "  -3.51%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.366"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -3.06%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "5.23"
$r.Style = "Normal"
$ws.Range("E41").Value = "  -3.36%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "1.80"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -4.43%  "
$ws.Range("E43").Value = "  -0.57%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "2.50"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -5.19%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "40.15"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -1.91%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "0.0₆0296"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -4.59%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "0.590"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -1.21%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "152.41"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -1.76%  "
$r = $ws.Range("D50")
$r.NumberFormat = "@"
$r.Value = "3.83"
$r.Style = "Normal"
$ws.Range("E50").Value = "  -2.30%  "
$ws.Range("E51").Value = "  -3.34%  "
